# Rotates the course-plan text blocks through the document (see commit diff):
# each paragraph/run receives the text that used to belong to the "next" block
# in the logical sequence, with the old "Docente" text wrapping around to
# "Bibliografia". Each Find is scoped to the exact paragraph Range so there is
# never any ambiguity about which run gets updated.
$d = $word.ActiveDocument

# Objetivos (paragraph 6): replace old block "B" with block "C"
$rng = $d.Paragraphs.Item(6).Range
$ok = $rng.Find.Execute('Aplicar as técnicas de modelação matemática no estudo de processos de tratamento de águas de abastecimento e residuárias. Fornecer ao aluno condições para uma análise matemática dos sistemas de tratamento de resíduos através de fundamentos de modelagem de fenômenos físicos e bioquímicos. Desenvolver a capacidade de uso de modelos matemáticos na simulação de processos empregados no tratamento de águas.', $true, $true, $false, $false, $false, $true, 1, $false, 'Modelagem matemática de processos de tratamento de águas residuárias: sedimentação, aeração, reatores aeróbios, reatores anaeróbios. Modelos matematicos de processos de tratamento de águas de abastecimento: floculação e filtração. Calibração e validação de modelos.', 2)
if (-not $ok) { throw "Replacement failed: Objetivos" }

# Docente (paragraph 8): replace old block "A" with block "B"
$rng = $d.Paragraphs.Item(8).Range
$ok = $rng.Find.Execute('229266 - Adriano Francisco Siqueira', $true, $true, $false, $false, $false, $true, 1, $false, 'Aplicar as técnicas de modelação matemática no estudo de processos de tratamento de águas de abastecimento e residuárias. Fornecer ao aluno condições para uma análise matemática dos sistemas de tratamento de resíduos através de fundamentos de modelagem de fenômenos físicos e bioquímicos. Desenvolver a capacidade de uso de modelos matemáticos na simulação de processos empregados no tratamento de águas.', 2)
if (-not $ok) { throw "Replacement failed: Docente" }

# Programa resumido (paragraph 10): replace old block "C" with block "D"
$rng = $d.Paragraphs.Item(10).Range
$ok = $rng.Find.Execute('Modelagem matemática de processos de tratamento de águas residuárias: sedimentação, aeração, reatores aeróbios, reatores anaeróbios. Modelos matematicos de processos de tratamento de águas de abastecimento: floculação e filtração. Calibração e validação de modelos.', $true, $true, $false, $false, $false, $true, 1, $false, '1- Dinâmica de processos físico-químicos e biológicos. 2- Revisão das equações fundamentais: cinética bio-química e conservação da massa. 3- Fundamentos dos modelos de floculação. 4- Modelos dinâmicos do processo de sedimentação. 5- Fundamentos dos modelos de filtração: a equação de Darcy e os modelos de resistência à filtração. 6- Modelos dinâmicos do processo de oxigenação de águas com e sem consumo simultâneo de oxigênio. 7- Fundamentos do modelo de tratamento de águas residuárias por lodos ativados. 8- Fundamentos dos modelos de digestão anaeróbia. 9- Calibração e validação de modelos.', 2)
if (-not $ok) { throw "Replacement failed: Programa resumido" }

# Programa (paragraph 12): replace old block "D" with block "E"
$rng = $d.Paragraphs.Item(12).Range
$ok = $rng.Find.Execute('1- Dinâmica de processos físico-químicos e biológicos. 2- Revisão das equações fundamentais: cinética bio-química e conservação da massa. 3- Fundamentos dos modelos de floculação. 4- Modelos dinâmicos do processo de sedimentação. 5- Fundamentos dos modelos de filtração: a equação de Darcy e os modelos de resistência à filtração. 6- Modelos dinâmicos do processo de oxigenação de águas com e sem consumo simultâneo de oxigênio. 7- Fundamentos do modelo de tratamento de águas residuárias por lodos ativados. 8- Fundamentos dos modelos de digestão anaeróbia. 9- Calibração e validação de modelos.', $true, $true, $false, $false, $false, $true, 1, $false, 'Aulas em sala de ensino informatizado, com auxílio de softwares para modelagem matemática.', 2)
if (-not $ok) { throw "Replacement failed: Programa" }

# Avaliacao - Norma de recuperacao run (paragraph 14): replace old block "G" with block "H"
$rng = $d.Paragraphs.Item(14).Range
$ok = $rng.Find.Execute('Prova única com nota igual ou superior a 5,0.', $true, $true, $false, $false, $false, $true, 1, $false, 'Pinto, José Carlos e Lage, Paulo Laranjeira C. Métodos Numéricos em Problemas de Engenharia Química. Rio de Janeiro, E-papers Serviços Editorias, 2001.Weber Jr., W. J. e DiGianno, F.A Process Dynamics in Environmental Systems.New York, J. Wiley & Sons. 1996.Garcia, Claudio. Modelagem e Simulação de Processos Industriais e de Sistemas Eletromecânicos. São Paulo, Edusp. 1997.Dochain, Denis e Vanrolleghem, Peter. A. Dynamical Modelling and Estimation in Wastewater Treatment Processes. London, IWA Publishing, 2001', 2)
if (-not $ok) { throw "Replacement failed: Avaliacao - Norma de recuperacao run" }

# Avaliacao - Criterio run (paragraph 14): replace old block "F" with block "G"
$rng = $d.Paragraphs.Item(14).Range
$ok = $rng.Find.Execute('A avaliação dos alunos deverá ser feita com base em exercícios resolvidos em casa e relatórios de atividades práticas.', $true, $true, $false, $false, $false, $true, 1, $false, 'Prova única com nota igual ou superior a 5,0.', 2)
if (-not $ok) { throw "Replacement failed: Avaliacao - Criterio run" }

# Avaliacao - Metodo run (paragraph 14): replace old block "E" with block "F"
$rng = $d.Paragraphs.Item(14).Range
$ok = $rng.Find.Execute('Aulas em sala de ensino informatizado, com auxílio de softwares para modelagem matemática.', $true, $true, $false, $false, $false, $true, 1, $false, 'A avaliação dos alunos deverá ser feita com base em exercícios resolvidos em casa e relatórios de atividades práticas.', 2)
if (-not $ok) { throw "Replacement failed: Avaliacao - Metodo run" }

# Bibliografia (paragraph 16): replace old block "H" with block "A"
$rng = $d.Paragraphs.Item(16).Range
$ok = $rng.Find.Execute('Pinto, José Carlos e Lage, Paulo Laranjeira C. Métodos Numéricos em Problemas de Engenharia Química. Rio de Janeiro, E-papers Serviços Editorias, 2001.Weber Jr., W. J. e DiGianno, F.A Process Dynamics in Environmental Systems.New York, J. Wiley & Sons. 1996.Garcia, Claudio. Modelagem e Simulação de Processos Industriais e de Sistemas Eletromecânicos. São Paulo, Edusp. 1997.Dochain, Denis e Vanrolleghem, Peter. A. Dynamical Modelling and Estimation in Wastewater Treatment Processes. London, IWA Publishing, 2001', $true, $true, $false, $false, $false, $true, 1, $false, '229266 - Adriano Francisco Siqueira', 2)
if (-not $ok) { throw "Replacement failed: Bibliografia" }
